# Edit script: reflects "change code create groups Whatsapp" commit.
# - Bulk-normalizes the I (horas inicio-bloque) / J (horas fin-bloque) columns
#   for the existing schedule rows (2-34) to the new group-size numbers,
#   re-applying the grey "highlight" formatting consistently down the range.
# - Adds two new course rows (35, 36: PROPAGACION DE PLANTAS / MECANIZACION
#   AGRICOLA) at the bottom of the table, including a live hyperlink on G35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cargas")

# ---------------------------------------------------------------------
# 1) Re-apply the existing I2:J2 formatting (grey highlight) across the
#    whole I2:J34 block so every row ends up visually consistent.
# ---------------------------------------------------------------------
$ws.Range("I2:J2").Copy()
$ws.Range("I2:J34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Update the group-size values in columns I/J for every existing row.
#    Row 2 and row 14 are special-cased; everything else becomes 5/16.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 34; $r++) {
    if ($r -eq 2) {
        $ws.Cells.Item($r, 9).Value2 = 14
        $ws.Cells.Item($r, 10).Value2 = 15
    } elseif ($r -eq 14) {
        $ws.Cells.Item($r, 9).Value2 = 16
        $ws.Cells.Item($r, 10).Value2 = 16
    } else {
        $ws.Cells.Item($r, 9).Value2 = 5
        $ws.Cells.Item($r, 10).Value2 = 16
    }
}

# ---------------------------------------------------------------------
# 3) Append the two new course rows (35 and 36).
# ---------------------------------------------------------------------

# --- Row 35: PROPAGACION DE PLANTAS ---
$ws.Range("A35").Value2 = "PROPAGACIÓN DE PLANTAS"
$ws.Range("B35").Value2 = "VI"
$ws.Range("C35").Value2 = "T"
$ws.Range("D35").Value2 = "A"
$ws.Range("E35").Value2 = 48
$ws.Range("F35").Value2 = "ALIAGA BARRERA ISAAC NOLBERTO"
$ws.Range("G35").Value2 = "https://sivireno.undc.edu.pe/index_home.php?s=asistencia_cursos.php&id_dcl=64"
$ws.Hyperlinks.Add($ws.Range("G35"), "https://sivireno.undc.edu.pe/index_home.php?s=asistencia_cursos.php&id_dcl=64") | Out-Null
$ws.Range("H35").Value2 = "SI"
$ws.Range("I35").Value2 = 5
$ws.Range("J35").Value2 = 16
$ws.Range("K34").Copy()
$ws.Range("K35").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K35").Value2 = "MIERCOLES"
$ws.Range("L35").Value2 = 0.65625
$ws.Range("M35").Value2 = 0.71875
$ws.Range("N34").Copy()
$ws.Range("N35").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N35").Value2 = "1-2"

# --- Row 36: MECANIZACION AGRICOLA ---
$ws.Range("A36").Value2 = "MECANIZACIÓN AGRÍCOLA"
$ws.Range("B36").Value2 = "V"
$ws.Range("C36").Value2 = "T"
$ws.Range("D36").Value2 = "A"
$ws.Range("E36").Value2 = 22
$ws.Range("F36").Value2 = "GARCIA RUIZ MARIA LUISA"
$ws.Range("G36").Value2 = "https://sivireno.undc.edu.pe/index_home.php?s=asistencia_cursos.php&id_dcl=63"
$h36 = $ws.Hyperlinks.Add($ws.Range("G36"), "https://sivireno.undc.edu.pe/index_home.php?s=asistencia_cursos.php&id_dcl=63")
$h36.Delete()
$ws.Range("H36").Value2 = "SI"
$ws.Range("I36").Value2 = 5
$ws.Range("J36").Value2 = 16
$ws.Range("K34").Copy()
$ws.Range("K36").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K36").Value2 = "MIERCOLES"
$ws.Range("L36").Value2 = 0.65625
$ws.Range("M36").Value2 = 0.71875
$ws.Range("N34").Copy()
$ws.Range("N36").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N36").Value2 = "1-2"

# ---------------------------------------------------------------------
# 4) Update the view: scroll so row ~22 is near the top and select L37,
#    matching where the author ended up after adding the new rows.
# ---------------------------------------------------------------------
$ws.Range("L37").Select()

Write-Host "edit complete"
